$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data
$ws.Range("A5").Value = "back to master"

# Update the selection to match the diff (A6)
$ws.Range("A6").Select()
